# Insert a new weekly price record for "Vega Monumental Concepción - Brócoli"
# at row 108 (pushing the existing rows 108-178 down to 109-179).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 108..178 down by one, creating a blank row 108.
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record.
$ws.Range("A108").Value = 11
$ws.Range("B108").Value = "Vega Monumental Concepción"
$ws.Range("C108").Value = "Bíobío"
$ws.Range("D108").Value = 44518
$ws.Range("E108").Value = 8
$ws.Range("F108").Value = 100112023
$ws.Range("G108").Value = "Brócoli"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 1100
$ws.Range("K108").Value = 700
$ws.Range("L108").Value = 800
$ws.Range("M108").Value = 745
$ws.Range("N108").Value = "$/unidad"
$ws.Range("O108").Value = "Región Metropolitana"
$ws.Range("P108").Value = 745
$ws.Range("Q108").Value = 1
$ws.Range("R108").Value = "Hortaliza"
